# Apply updated cryptos data (Price and Volume(1h) columns) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    # Assigns a string value to a cell while forcing text interpretation so that
    # numeric-looking strings (e.g. "218.38") are not auto-converted to numbers,
    # matching the inline-string cells produced by the source data feed.
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "26.285.77"
$ws.Range("E2").Value = "  -5.23%  "
Set-TextCell $ws.Range("D3") "1.675.35"
$ws.Range("E3").Value = "  -2.70%  "
$ws.Range("E4").Value = "  +0.33%  "
Set-TextCell $ws.Range("D5") "218.38"
$ws.Range("E5").Value = "  -2.65%  "
Set-TextCell $ws.Range("D6") "0.5140"
$ws.Range("E6").Value = "  -9.39%  "
$ws.Range("E7").Value = "  +0.31%  "
Set-TextCell $ws.Range("D8") "0.2667"
$ws.Range("E8").Value = "  -1.30%  "
Set-TextCell $ws.Range("D9") "0.06416"
$ws.Range("E9").Value = "  -1.97%  "
Set-TextCell $ws.Range("D10") "21.57"
$ws.Range("E10").Value = "  -5.04%  "
Set-TextCell $ws.Range("D11") "0.07374"
$ws.Range("E11").Value = "  -1.74%  "
Set-TextCell $ws.Range("D12") "1.674.97"
$ws.Range("E12").Value = "  -2.71%  "
Set-TextCell $ws.Range("D13") "4.570"
$ws.Range("E13").Value = "  -1.92%  "
Set-TextCell $ws.Range("D14") "0.5841"
$ws.Range("E14").Value = "  -1.57%  "
Set-TextCell $ws.Range("D15") "1.902.41"
$ws.Range("E15").Value = "  -2.76%  "
Set-TextCell $ws.Range("D16") "0.000008699"
$ws.Range("E16").Value = "  +2.16%  "
Set-TextCell $ws.Range("D17") "65.02"
$ws.Range("E17").Value = "  -11.83%  "
Set-TextCell $ws.Range("D18") "26.372.16"
$ws.Range("E18").Value = "  -4.85%  "
Set-TextCell $ws.Range("D19") "4.975"
$ws.Range("E19").Value = "  -5.34%  "
$ws.Range("E20").Value = "  +0.11%  "
Set-TextCell $ws.Range("D21") "10.86"
$ws.Range("E21").Value = "  -2.86%  "
Set-TextCell $ws.Range("D22") "190.87"
$ws.Range("E22").Value = "  -4.62%  "
Set-TextCell $ws.Range("D23") "6.257"
$ws.Range("E23").Value = "  -3.88%  "
Set-TextCell $ws.Range("D25") "144.25"
$ws.Range("E25").Value = "  -3.41%  "
Set-TextCell $ws.Range("D26") "7.695"
$ws.Range("E26").Value = "  -2.88%  "
Set-TextCell $ws.Range("D27") "0.1185"
$ws.Range("E27").Value = "  -2.20%  "
Set-TextCell $ws.Range("D28") "15.71"
$ws.Range("E28").Value = "  -1.93%  "
Set-TextCell $ws.Range("D29") "0.05897"
$ws.Range("E29").Value = "  -3.70%  "
Set-TextCell $ws.Range("D30") "1.274"
$ws.Range("E30").Value = "  -6.50%  "
Set-TextCell $ws.Range("D31") "1.327"
$ws.Range("E31").Value = "  -3.93%  "
Set-TextCell $ws.Range("D32") "3.537"
$ws.Range("E32").Value = "  -3.98%  "
Set-TextCell $ws.Range("D33") "3.524"
$ws.Range("E33").Value = "  -4.37%  "
Set-TextCell $ws.Range("D34") "1.649"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  -1.08%  "
Set-TextCell $ws.Range("D36") "0.6037"
$ws.Range("E36").Value = "  -5.87%  "
Set-TextCell $ws.Range("D37") "2.359"
$ws.Range("E37").Value = "  -2.54%  "
Set-TextCell $ws.Range("D38") "2.649"
$ws.Range("E38").Value = "  -1.22%  "
Set-TextCell $ws.Range("D39") "0.01623"
$ws.Range("E39").Value = "  -2.44%  "
Set-TextCell $ws.Range("D40") "6.051"
$ws.Range("E40").Value = "  -1.18%  "
Set-TextCell $ws.Range("D41") "1.080.81"
$ws.Range("E41").Value = "  -2.90%  "
Set-TextCell $ws.Range("D42") "0.8694"
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("E43").Value = "  +0.63%  "
Set-TextCell $ws.Range("D44") "99.86"
$ws.Range("E44").Value = "  +0.73%  "
Set-TextCell $ws.Range("D45") "1.823.71"
$ws.Range("E45").Value = "  -2.48%  "
Set-TextCell $ws.Range("D46") "0.00000000112"
$ws.Range("E46").Value = "  +5.47%  "
Set-TextCell $ws.Range("D47") "56.19"
$ws.Range("E47").Value = "  -4.29%  "
Set-TextCell $ws.Range("D48") "1.010"
$ws.Range("E48").Value = "  +1.22%  "
Set-TextCell $ws.Range("D49") "8.119"
$ws.Range("E49").Value = "  -0.84%  "
Set-TextCell $ws.Range("D50") "0.4299"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("E51").Value = "  -3.10%  "
